$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db_table")

# 1. Batch rename exp001 -> exp002 for the duplicate-cell rows (10-16), column C.
#    Doing this first so "exp002" becomes shared-string index 121 (matches target ordering).
foreach ($r in 10..16) {
    $ws.Range("C" + $r).Value = "exp002"
}

# 2. N3: set to existing "test" string.
$ws.Range("N3").Value = "test"

# 3. AV3: set to new "test comment general" string (becomes shared-string index 122).
$ws.Range("AV3").Value = "test comment general"

# 4. P4: set to text "1" (column already text-formatted, matches O4).
$ws.Range("P4").Value = "1"

# 5. T4: numeric value changed from 1 to 0.
$ws.Range("T4").Value = 0

# 6. AJ4: tiny float recalculation.
$ws.Range("AJ4").Value = 0.45438760079739898

# 7. New row 18: a "column index" helper row with an incrementing formula chain
#    from C18 (=1+B18) through BB18, plus summary cells at BC18:BE18.
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = "1"

$ws.Range("C18").Formula = "=1+B18"

$ws.Range("D18:BB18").FormulaR1C1 = "=1+RC[-1]"

$ws.Range("BC18").Value = 27
$ws.Range("BD18").Value = "28"
$ws.Range("BE18").Value = 28

# Apply the "index row" style (text numeric format, bold-less font, explicit
# fill applied) to C18:BB18, cloned from an existing cell using that
# combination (font without explicit color + applyFill), then force text
# number format on top so Excel reuses/creates the s="28" xf.
$ws.Range("AV3").Copy() | Out-Null
$ws.Range("C18:BB18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18:BB18").NumberFormat = "@"

# A18/BC18/BE18 use the plain bold-ish header style (s=6, same as AV3 before
# the text-format tweak); BD18 uses the default left column style (s=2).
$ws.Range("A18").Style = $ws.Range("T4").Style
$ws.Range("BC18").Style = $ws.Range("T4").Style
$ws.Range("BE18").Style = $ws.Range("T4").Style
$ws.Range("BD18").Style = $ws.Range("B18").Style

# 8. Column width adjustments.
$ws.Columns.Item(3).ColumnWidth = 21.7109375
$ws.Columns.Item(16).ColumnWidth = 16.5703125
$ws.Columns.Item(40).ColumnWidth = 22
$ws.Columns.Item(42).ColumnWidth = 19.7109375
$ws.Columns.Item(43).ColumnWidth = 26.140625
$ws.Columns.Item(44).ColumnWidth = 18.28515625
$ws.Range("AS1:AU1").EntireColumn.Hidden = $false

# 9. Sheet view / selection changes on db_table.
$av = $ws.Application.ActiveWindow
$av.ScrollColumn = $ws.Range("AK3").Column
$ws.Range("AN4:AN6").Select()

# 10. Notes sheet view changes.
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A46").Select()
$notesWin = $notes.Application.ActiveWindow
$notesWin.ScrollRow = 13

Write-Output "edit applied"
